{"js": "// Translate specific English strings in the document to French, matching\n// the OOXML diff exactly (whole-run text replacements).\n\nconst replacements = [\n  {\n    find: \"A Lesson on Conway\\u2019s Game of Life [Samuel Okoth]\",\n    replace: \"Une le\\u00e7on sur le jeu de la vie de Conway [Samuel Flot]\",\n  },\n  {\n    find:\n      \"Define an interesting and unpredictable cell automaton. For example, discover some configurations that last for a long time before dying and other configurations to go on forever without allowing cycles.\",\n    replace:\n      \"D\\u00e9finir un automate cellulaire int\\u00e9ressant et impr\\u00e9visible. Par exemple, d\\u00e9couvrez des configurations qui durent longtemps avant de mourir et d'autres configurations qui durent \\u00e9ternellement sans permettre de cycles.\",\n  },\n  {\n    find: \"(Leave \",\n    replace: \"(Laissez \",\n  },\n  {\n    find: \" for the facilitators that will use it)\",\n    replace: \" pour les facilitateurs qui l'utiliseront)\",\n  },\n  {\n    find: \"Paper to draw square grids, 2 different coloured post-its\",\n    replace:\n      \"Papier pour dessiner des grilles carr\\u00e9es, 2 post-it de couleurs diff\\u00e9rentes\",\n  },\n  {\n    find: \"Introduction of the second experiment\",\n    replace: \"Introduction de la deuxi\\u00e8me exp\\u00e9rience\",\n  },\n];\n\nfor (const { find, replace } of replacements) {\n  const results = context.document.body.search(find, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items,text\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Translate specific English strings to French, matching the OOXML diff\n# exactly (whole-run text replacements). Direct Range.Text assignment is\n# used (instead of Find.Execute's Replacement.Text) so Word's smart-quote\n# autocorrect does not mangle the straight apostrophes in the French text.\n\n$d = $word.ActiveDocument\n\nfunction Replace-AllText($findText, $replaceText) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.MatchCase = $true\n    $rng.Find.MatchWholeWord = $false\n    while ($rng.Find.Execute($findText)) {\n        $rng.Text = $replaceText\n        $rng.Collapse(0)\n        $rng.End = $d.Content.End\n    }\n}\n\nReplace-AllText \"A Lesson on Conway\u2019s Game of Life [Samuel Okoth]\" \"Une le\u00e7on sur le jeu de la vie de Conway [Samuel Flot]\"\n\nReplace-AllText \"Define an interesting and unpredictable cell automaton. For example, discover some configurations that last for a long time before dying and other configurations to go on forever without allowing cycles.\" \"D\u00e9finir un automate cellulaire int\u00e9ressant et impr\u00e9visible. Par exemple, d\u00e9couvrez des configurations qui durent longtemps avant de mourir et d'autres configurations qui durent \u00e9ternellement sans permettre de cycles.\"\n\nReplace-AllText \"(Leave \" \"(Laissez \"\n\nReplace-AllText \" for the facilitators that will use it)\" \" pour les facilitateurs qui l'utiliseront)\"\n\nReplace-AllText \"Paper to draw square grids, 2 different coloured post-its\" \"Papier pour dessiner des grilles carr\u00e9es, 2 post-it de couleurs diff\u00e9rentes\"\n\nReplace-AllText \"Introduction of the second experiment\" \"Introduction de la deuxi\u00e8me exp\u00e9rience\"\n"}
